$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename header row: "<name>_old" -> "<name>_FV2404", "<name>_new" -> "<name>_FV2410"
#    (column K stays "diff"). These ten base names live in columns A-J (the
#    "_old" / FV2404 side) and L-U (the "_new" / FV2410 side).
# ---------------------------------------------------------------------------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row (row 1) so it stays pinned while scrolling.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3) Turn the used range A1:U64 into an Excel Table ("Table1").
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$table.Name = "Table1"
